$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("B44").Value = "RocketPoolETH"

# --- Column C (Link) updates ---
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"

# --- Column D (Price) updates: force text format to preserve exact digits ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.903.61"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.16"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.57"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6366"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3008"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07482"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.58"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.53"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.054"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6903"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.32"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009397"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.090"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.865.69"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.121.63"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.48"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.369"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.35"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1421"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.591"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06090"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.144"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.145"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.878"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.165"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.618"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.865"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01797"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.222.65"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9269"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.312"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.031.58"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.43"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5091"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.290"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4094"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1145"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +4.43%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +7.80%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("E14").Value = "  +4.87%  "
$ws.Range("E15").Value = "  +3.59%  "
$ws.Range("E16").Value = "  +4.66%  "
$ws.Range("E17").Value = "  +4.62%  "
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("E27").Value = "  +2.18%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("E30").Value = "  +9.82%  "
$ws.Range("E31").Value = "  +5.19%  "
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("E34").Value = "  +3.82%  "
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("E51").Value = "  +3.39%  "
